$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ID" row (row 2) -- Property/Id/建筑ID row is no longer used;
# remaining rows shift up automatically.
$ws.Rows.Item(2).Delete()

# Fix the Type column (B) for rows that actually hold string data instead
# of int (Prefab, NormalStateFunc, UpStateFunc, Desc -- now rows 4-7 after
# the delete above).
$ws.Range("B4").Value = "string"
$ws.Range("B5").Value = "string"
$ws.Range("B6").Value = "string"
$ws.Range("B7").Value = "string"

# Re-apply the TRUE/FALSE list validation for the rows below the table
# (now starting at row 8 instead of row 9).
$ws.Range("F8:F1048576").Validation.Delete()
$ws.Range("F8:F1048576").Validation.Add(3, 1, 1, """TRUE,FALSE""")

# Match the saved selection/active cell recorded in the workbook.
$ws.Range("G14").Select()
